$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("effort")

# Keep the date column formatting consistent with the rest of the table by
# copying the style of the preceding row's date cell
$ws.Range("A51").Copy() | Out-Null
$ws.Range("A52").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# New data row appended to the effort log table
$ws.Range("A52").Value = 41246
$ws.Range("B52").Value = 2.5
$ws.Range("D52").Value = 'Manual: New section "Task switches" started'

# Match the selection state recorded after the edit
$ws.Range("E52").Select()
